$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = 17
$ws.Range("D3").Value = 17
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 22
$ws.Range("C6").Value = 36
$ws.Range("D6").Value = 36
$ws.Range("C8").Value = 116
$ws.Range("D8").Value = 116
$ws.Range("C10").Value = 220
$ws.Range("D10").Value = 220
$ws.Range("C12").Value = 62
$ws.Range("D12").Value = 62
$ws.Range("C14").Value = 29
$ws.Range("D14").Value = 29
$ws.Range("C16").Value = 147
$ws.Range("D16").Value = 147
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 10
$ws.Range("C20").Value = 52
$ws.Range("D20").Value = 52
$ws.Range("C22").Value = 197
$ws.Range("D22").Value = 197
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 5
$ws.Range("C26").Value = 146
$ws.Range("D26").Value = 146
$ws.Range("C28").Value = 33
$ws.Range("D28").Value = 33
$ws.Range("C29").Value = 110
$ws.Range("D29").Value = 110
$ws.Range("C31").Value = 23
$ws.Range("D31").Value = 23
$ws.Range("C33").Value = 69
$ws.Range("D33").Value = 69
$ws.Range("C35").Value = 164
$ws.Range("D35").Value = 164
$ws.Range("C37").Value = 12
$ws.Range("D37").Value = 12
$ws.Range("C39").Value = 167
$ws.Range("D39").Value = 167
$ws.Range("C41").Value = 118
$ws.Range("D41").Value = 118
$ws.Range("C44").Value = 27
$ws.Range("D44").Value = 130
$ws.Range("C45").Value = 44
$ws.Range("D45").Value = 44
$ws.Range("C46").Value = 228
$ws.Range("D46").Value = 228
$ws.Range("C49").Value = 119
$ws.Range("D49").Value = 119
$ws.Range("C50").Value = 112
$ws.Range("D50").Value = 112
$ws.Range("C52").Value = 208
$ws.Range("D52").Value = 208
$ws.Range("C54").Value = 11
$ws.Range("D54").Value = 11
$ws.Range("C56").Value = 175
$ws.Range("D56").Value = 175
$ws.Range("C58").Value = 16
$ws.Range("D58").Value = 16
$ws.Range("C61").Value = 39
$ws.Range("D61").Value = 39
$ws.Range("C64").Value = 150
$ws.Range("D64").Value = 150
$ws.Range("C66").Value = 71
$ws.Range("D66").Value = 71
$ws.Range("C70").Value = 72
$ws.Range("D70").Value = 72
$ws.Range("C74").Value = 9
$ws.Range("D74").Value = 9
$ws.Range("C75").Value = 145
$ws.Range("D75").Value = 145
$ws.Range("C77").Value = 123
$ws.Range("D77").Value = 123
$ws.Range("C78").Value = 84.59999999999999
